# Update BOM (all components in stock on LCSC on 2/6/2020)
# Row 3 corresponds to the 309K / R2,R4 / R0402 line item. Its manufacturer
# part number, manufacturer, and supplier part number are being refreshed to
# the currently-in-stock LCSC part.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM_TSDZ2-ESP32-v3")

$ws.Range("F3").Value = "WR04X3093FTL"
$ws.Range("G3").Value = "Walsin Tech Corp"
$ws.Range("I3").Value = "C334683"

# Reflect the author's final on-screen selection: the whole row 3 selected.
$ws.Rows.Item(3).Select()
